$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: C36 currently holds the CNPJ as text "04252502000160".
# Change it to the numeric value 4252502000160 (no leading zero).
$ws.Cells.Item(36, 3).Value = 4252502000160

# Insert a new row 37 with a new log entry (duplicate-looking send later the same day).
$ws.Cells.Item(37, 1).Value = "04/07/2025 17:26:42"
$ws.Cells.Item(37, 2).Value = "Ima Industria"

# C37 keeps the CNPJ as text, preserving the leading zero.
$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "04252502000160"

$ws.Cells.Item(37, 4).Value = "denissonfhsilva@gmail.com"
$ws.Cells.Item(37, 5).Value = "893-ExtratoMensal-052025.pdf"
